$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("F10").Value = 43.54
$ws1.Range("F22").Value = "1 de 20"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F10").Value = -329.51
$ws2.Range("F22").Value = -329.51

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# Excel's ColumnWidth property adds ~0.8333 (5/6) chars of padding when stored
# in the OOXML "width" attribute, so subtract it to land exactly on 23.
$ws3.Range("E:E").ColumnWidth = 22.166666666666668
$ws3.Range("D6").Value = 43.54
$ws3.Range("E6").Value = 63.27999999999999
$ws3.Range("F6").Value = 0.4076015727391875
$ws3.Range("D19").Value = -329.51
$ws3.Range("E19").Value = 50716.70762291769
$ws3.Range("F19").Value = -0.006539557973951076
